# Apply updated enrollment counts to the "Inscricoes" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 12
$ws.Range("E12").Value = 409
$ws.Range("F12").Value = 222
$ws.Range("H12").Value = 222

# Row 16
$ws.Range("E16").Value = 173

# Row 17
$ws.Range("E17").Value = 80
$ws.Range("F17").Value = 38
$ws.Range("H17").Value = 38

# Row 20
$ws.Range("E20").Value = 80

# Row 26
$ws.Range("E26").Value = 124
$ws.Range("F26").Value = 76
$ws.Range("H26").Value = 76

# Row 30
$ws.Range("E30").Value = 182
$ws.Range("F30").Value = 101
$ws.Range("H30").Value = 101

# Row 43
$ws.Range("F43").Value = 52
$ws.Range("H43").Value = 52

# Row 48
$ws.Range("E48").Value = 181

# Row 51
$ws.Range("E51").Value = 202
